$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 429.9
$ws.Range("I12").Value = 433.22223
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 433.22223
$ws.Range("L12").Value = 400
$ws.Range("M12").Value = -263.22223
$ws.Range("N12").Value = -740

$ws.Range("H62").Value = 1766.5555
$ws.Range("I62").Value = 1649.875
$ws.Range("K62").Value = 1649.875
$ws.Range("M62").Value = -1025.875

$ws.Range("H65").Value = 1766.5555
$ws.Range("I65").Value = 1649.875
$ws.Range("K65").Value = 8249.375
$ws.Range("M65").Value = -5129.375

$ws.Range("H92").Value = 665.94116
$ws.Range("I92").Value = 573.4167
$ws.Range("J92").Value = 888
$ws.Range("K92").Value = 573.4167
$ws.Range("L92").Value = 888
$ws.Range("M92").Value = 674.5833
$ws.Range("N92").Value = -3384

$ws.Range("H96").Value = 1115.3572
$ws.Range("I96").Value = 1261.7
$ws.Range("J96").Value = 749.5
$ws.Range("K96").Value = 3785.1
$ws.Range("L96").Value = 2248.5
$ws.Range("M96").Value = -2412.1
$ws.Range("N96").Value = -4994.5

$ws.Range("H116").Value = 1936170.8
$ws.Range("I116").Value = 7144612
$ws.Range("K116").Value = 7144612
$ws.Range("M116").Value = -7141170

$ws.Range("H131").Value = 1865

$ws.Range("H132").Value = 3954839
$ws.Range("I132").Value = 4330538
$ws.Range("K132").Value = 12991614
$ws.Range("M132").Value = -12989084

$ws.Range("H141").Value = 2748.7568
$ws.Range("I141").Value = 1386.0377
$ws.Range("J141").Value = 6188
$ws.Range("K141").Value = 4158.1131
$ws.Range("L141").Value = 18564
$ws.Range("M141").Value = 1021.8869
$ws.Range("N141").Value = -28924

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1380.5358
$ws.Range("I61").Value = 1216.55
$ws.Range("J61").Value = 1790.5
$ws.Range("K61").Value = 1216.55
$ws.Range("L61").Value = 1790.5
$ws.Range("M61").Value = -1004.55
$ws.Range("N61").Value = -2214.5

$ws.Range("H74").Value = 1337.6666
$ws.Range("I74").Value = 1433.6296
$ws.Range("J74").Value = 1121.75
$ws.Range("K74").Value = 1433.6296
$ws.Range("L74").Value = 1121.75
$ws.Range("M74").Value = -559.6296
$ws.Range("N74").Value = -2869.75

$ws.Range("H77").Value = 1337.6666
$ws.Range("I77").Value = 1433.6296
$ws.Range("J77").Value = 1121.75
$ws.Range("K77").Value = 7168.148
$ws.Range("L77").Value = 5608.75
$ws.Range("M77").Value = -2800.148
$ws.Range("N77").Value = -14344.75

$ws.Range("H102").Value = 1561.3334
$ws.Range("I102").Value = 1561.3334
$ws.Range("K102").Value = 1561.3334
$ws.Range("M102").Value = 60.66660000000002

$ws.Range("H136").Value = 1380.5358
$ws.Range("I136").Value = 1216.55
$ws.Range("J136").Value = 1790.5
$ws.Range("K136").Value = 3649.65
$ws.Range("L136").Value = 5371.5
$ws.Range("M136").Value = -1099.65
$ws.Range("N136").Value = -10471.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1625208.8
$ws.Range("I105").Value = 3789520.2
$ws.Range("J105").Value = 1975
$ws.Range("K105").Value = 3789520.2
$ws.Range("L105").Value = 1975
$ws.Range("M105").Value = -3787773.2
$ws.Range("N105").Value = -5469

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2190.4167
$ws.Range("I132").Value = 2062.1428
$ws.Range("J132").Value = 2370
$ws.Range("K132").Value = 6186.428400000001
$ws.Range("L132").Value = 7110
$ws.Range("M132").Value = -3656.428400000001
$ws.Range("N132").Value = -12170

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 354.02856
$ws.Range("I5").Value = 235.39285
$ws.Range("J5").Value = 828.5714
$ws.Range("K5").Value = 706.1785500000001
$ws.Range("L5").Value = 2485.7142
$ws.Range("M5").Value = -594.1785500000001
$ws.Range("N5").Value = -2709.7142

$ws.Range("H131").Value = 878.383
$ws.Range("I131").Value = 453.92307
$ws.Range("J131").Value = 1040.6765
$ws.Range("K131").Value = 1361.76921
$ws.Range("L131").Value = 3122.0295
$ws.Range("M131").Value = 3678.23079
$ws.Range("N131").Value = -13202.0295

$ws.Range("H135").Value = 354.02856
$ws.Range("I135").Value = 235.39285
$ws.Range("J135").Value = 828.5714
$ws.Range("K135").Value = 2118.53565
$ws.Range("L135").Value = 7457.1426
$ws.Range("M135").Value = 416.4643499999997
$ws.Range("N135").Value = -12527.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 45000
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

$ws.Range("H85").Value = 45000
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 377.96155
$ws.Range("I22").Value = 346.91666
$ws.Range("J22").Value = 750.5
$ws.Range("K22").Value = 346.91666
$ws.Range("L22").Value = 750.5
$ws.Range("M22").Value = -51.91665999999998
$ws.Range("N22").Value = -1340.5

$ws.Range("H27").Value = 377.96155
$ws.Range("I27").Value = 346.91666
$ws.Range("J27").Value = 750.5
$ws.Range("K27").Value = 346.91666
$ws.Range("L27").Value = 750.5
$ws.Range("M27").Value = -239.91666
$ws.Range("N27").Value = -964.5

$ws.Range("H93").Value = 11217.818
$ws.Range("I93").Value = 17883.666
$ws.Range("K93").Value = 17883.666
$ws.Range("M93").Value = -16635.666

$ws.Range("H132").Value = 1671.725
$ws.Range("I132").Value = 1030.1562
$ws.Range("J132").Value = 4238
$ws.Range("K132").Value = 3090.4686
$ws.Range("L132").Value = 12714
$ws.Range("M132").Value = -560.4685999999997
$ws.Range("N132").Value = -17774

$ws.Range("H136").Value = 1997.7567
$ws.Range("I136").Value = 1214.2307
$ws.Range("J136").Value = 2422.1667
$ws.Range("K136").Value = 3642.6921
$ws.Range("L136").Value = 7266.500100000001
$ws.Range("M136").Value = -1092.6921
$ws.Range("N136").Value = -12366.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 63580
$ws.Range("J16").Value = 63580
$ws.Range("L16").Value = 63580
$ws.Range("N16").Value = -64164

$ws.Range("H80").Value = 90060.2
$ws.Range("J80").Value = 90060.2
$ws.Range("L80").Value = 90060.2
$ws.Range("N80").Value = -92056.2

$ws.Range("H83").Value = 90060.2
$ws.Range("J83").Value = 90060.2
$ws.Range("L83").Value = 270180.6
$ws.Range("N83").Value = -280164.6

$ws.Range("I96").Value = 1903
$ws.Range("J96").Value = 2666.6667
$ws.Range("K96").Value = 1903
$ws.Range("L96").Value = 2666.6667
$ws.Range("M96").Value = -530
$ws.Range("N96").Value = -5412.6667
